# Update "paises" (countries) COVID stats sheet + refresh the "datos
# actualizados" timestamp, per the later data pull (01:26 -> 02:43).
#
# The underlying table (sheet "Pais") is ranked descending by total
# cases (column B). Between the two snapshots some countries' totals
# changed enough to re-order a few neighbouring rows; for those rows
# both the country name (col A) and the stats (cols B-H) are rewritten
# to reflect the new rank. For rows whose rank didn't change, only the
# stats are refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 22 de Octubre de 2020 a las 02:43"

# --- Rows whose rank is unchanged: stats-only refresh ------------------
# Row 4  -> Estados Unidos
$ws.Range("B4").Value = 8581997
$ws.Range("C4").Value = 60841
$ws.Range("D4").Value = 5586588
$ws.Range("E4").Value = 2768052
$ws.Range("G4").Value = 1173
$ws.Range("H4").Value = 227357

# Row 6  -> Brasil
$ws.Range("B6").Value = 5300649
$ws.Range("C6").Value = 25832
$ws.Range("D6").Value = 4756489
$ws.Range("E6").Value = 388701
$ws.Range("G6").Value = 571
$ws.Range("H6").Value = 155459

# Row 32 -> Canada
$ws.Range("B32").Value = 205954
$ws.Range("C32").Value = 2266
$ws.Range("D32").Value = 173514
$ws.Range("E32").Value = 22614
$ws.Range("G32").Value = 32
$ws.Range("H32").Value = 9826

# Row 150 -> Gambia
$ws.Range("B150").Value = 3657
$ws.Range("C150").Value = 2
$ws.Range("E150").Value = 881

# Row 156 -> Benin
$ws.Range("B156").Value = 2557
$ws.Range("C156").Value = 61
$ws.Range("E156").Value = 186

# Row 162 -> Yemen
$ws.Range("D162").Value = 1344
$ws.Range("E162").Value = 116

# Row 169 -> Santo Tome y Principe
$ws.Range("B169").Value = 935
$ws.Range("C169").Value = 2
$ws.Range("E169").Value = 22

# --- Rows 116-118: Angola overtakes Lituania & Zimbabue -----------------
$ws.Range("A116").Value = "Angola"
$ws.Range("B116").Value = 8338
$ws.Range("C116").Value = 289
$ws.Range("D116").Value = 3040
$ws.Range("E116").Value = 5043
$ws.Range("G116").Value = 4
$ws.Range("H116").Value = 255

$ws.Range("A117").Value = "Lituania"
$ws.Range("B117").Value = 8239
$ws.Range("C117").Value = 311
$ws.Range("D117").Value = 3599
$ws.Range("E117").Value = 4520
$ws.Range("G117").Value = 2
$ws.Range("H117").Value = 120

$ws.Range("A118").Value = "Zimbabue"
$ws.Range("B118").Value = 8215
$ws.Range("C118").Value = 28
$ws.Range("D118").Value = 7725
$ws.Range("E118").Value = 254
$ws.Range("G118").Value = 3
$ws.Range("H118").Value = 236

# --- Rows 133-140: Polinesia Francesa overtakes Congo..Rep. Africa Central
$ws.Range("A133").Value = "Polinesia Francesa"
$ws.Range("B133").Value = 5161
$ws.Range("C133").Value = 551
$ws.Range("D133").Value = 3536
$ws.Range("E133").Value = 1606
$ws.Range("G133").Value = 3
$ws.Range("H133").Value = 19

$ws.Range("A134").Value = "Congo"
$ws.Range("B134").Value = 5156
$ws.Range("D134").Value = 3887
$ws.Range("E134").Value = 1177
$ws.Range("H134").Value = 92

$ws.Range("A135").Value = "Surinam"
$ws.Range("B135").Value = 5150
$ws.Range("C135").Value = 6
$ws.Range("D135").Value = 4991
$ws.Range("E135").Value = 50
$ws.Range("H135").Value = 109

$ws.Range("A136").Value = "Guinea Ecuatorial"
$ws.Range("B136").Value = 5074
$ws.Range("C136").Value = 0
$ws.Range("D136").Value = 4954
$ws.Range("E136").Value = 37
$ws.Range("H136").Value = 83

$ws.Range("A137").Value = "Malta"
$ws.Range("B137").Value = 5026
$ws.Range("C137").Value = 155
$ws.Range("D137").Value = 3331
$ws.Range("E137").Value = 1649
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 46

$ws.Range("A138").Value = "Reunion"
$ws.Range("B138").Value = 5015
$ws.Range("C138").Value = 94
$ws.Range("D138").Value = 4445
$ws.Range("E138").Value = 551
$ws.Range("G138").Value = 2
$ws.Range("H138").Value = 19

$ws.Range("A139").Value = "Ruanda"
$ws.Range("B139").Value = 5012
$ws.Range("C139").Value = 16
$ws.Range("D139").Value = 4798
$ws.Range("E139").Value = 180
$ws.Range("H139").Value = 34

$ws.Range("A140").Value = "Republica de Africa Central"
$ws.Range("B140").Value = 4858
$ws.Range("D140").Value = 1924
$ws.Range("E140").Value = 2872
$ws.Range("H140").Value = 62
